$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow writes, then restore protection
$ws.Unprotect()

# Update the confidential disclosure date string (shared string used in A41)
$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-06 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for holdings rows 2-38
$ws.Range("D2").Value = 0.03181883764230603
$ws.Range("E2").Value = 0.002591121091725546
$ws.Range("D3").Value = 0.02866133447449705
$ws.Range("E3").Value = -0.006328463703795029
$ws.Range("D4").Value = 0.02902841752851028
$ws.Range("E4").Value = 0.0007709109597842101
$ws.Range("D5").Value = 0.06420105943702807
$ws.Range("E5").Value = 0.01095537739945085
$ws.Range("D6").Value = 0.01596752394628024
$ws.Range("E6").Value = 0.004130707383639409
$ws.Range("D7").Value = 0.0158591657399619
$ws.Range("E7").Value = 0.01002599331600429
$ws.Range("D8").Value = 0.02936860733204232
$ws.Range("E8").Value = 0.007820332865450252
$ws.Range("D9").Value = 0.0345764754726685
$ws.Range("E9").Value = 0.002952197115930533
$ws.Range("D10").Value = 0.02921431466869773
$ws.Range("E10").Value = 0.006994839542009412
$ws.Range("D11").Value = 0.03116201416705027
$ws.Range("E11").Value = -0.001587441573331128
$ws.Range("D12").Value = 0.01131283230095314
$ws.Range("E12").Value = -0.009717161200763491
$ws.Range("D13").Value = 0.01411208596417701
$ws.Range("E13").Value = 0.02921129503407993
$ws.Range("D14").Value = 0.01476478711636631
$ws.Range("E14").Value = -0.03536528617961832
$ws.Range("D15").Value = 0.009129571484517242
$ws.Range("E15").Value = 0.01255697944439671
$ws.Range("D16").Value = 0.007913093577715115
$ws.Range("E16").Value = 0.002455905336012476
$ws.Range("D17").Value = 0.0295982796171736
$ws.Range("E17").Value = 0.01230269266480977
$ws.Range("D18").Value = 0.02573468139840977
$ws.Range("E18").Value = 0.02106058063433447
$ws.Range("D19").Value = 0.03195094828153111
$ws.Range("E19").Value = 0.01539028660952901
$ws.Range("D20").Value = 0.03091938600942441
$ws.Range("E20").Value = 0.01587200812646805
$ws.Range("D21").Value = 0.04626306506498057
$ws.Range("E21").Value = 0.01044239076011788
$ws.Range("D22").Value = 0.03510059940540397
$ws.Range("E22").Value = 0.02334880599519051
$ws.Range("D23").Value = 0.03259088986232422
$ws.Range("E23").Value = 0.01665411835566921
$ws.Range("D24").Value = 0.03067891716388099
$ws.Range("E24").Value = 0.007198387561186292
$ws.Range("D25").Value = 0.01466369205068887
$ws.Range("E25").Value = -0.02771084337349383
$ws.Range("D26").Value = 0.01505354603211682
$ws.Range("E26").Value = -0.03935529301306617
$ws.Range("D27").Value = 0.03135753006105946
$ws.Range("E27").Value = 0.022035532295827
$ws.Range("D28").Value = 0.02969328934880055
$ws.Range("E28").Value = -0.02110217897186384
$ws.Range("D29").Value = 0.02902939903400229
$ws.Range("E29").Value = 0.01322676187771332
$ws.Range("D30").Value = 0.02854826504181704
$ws.Range("E30").Value = 0.009681567203690955
$ws.Range("D31").Value = 0.03405863317508193
$ws.Range("E31").Value = 0.004461043676729837
$ws.Range("D32").Value = 0.03107034155409616
$ws.Range("E32").Value = 0.0115713392174579
$ws.Range("D33").Value = 0.029443986953829
$ws.Range("E33").Value = 0.01404056162246481
$ws.Range("D34").Value = 0.03238968123646136
$ws.Range("E34").Value = 0.005987878787878875
$ws.Range("D35").Value = 0.03048634578634785
$ws.Range("E35").Value = 0
$ws.Range("D36").Value = 0.03092134902040844
$ws.Range("E36").Value = 0.02012442864398167
$ws.Range("D37").Value = 0.03335705304939034
$ws.Range("E37").Value = 0.002777647003436856
$ws.Range("D38").Value = 0.9999999999999999
$ws.Range("E38").Value = 0.006767185915788732

# Restore sheet protection
$ws.Protect()
